$d = $word.ActiveDocument

# 1. "Front-end" heading -> "Maxim (Front-end)"
$d.Content.Find.Execute("Front-end", $true, $false, $false, $false, $false, $true, 1, $false, "Maxim (Front-end)", 2)

# 2. "Backend" heading -> "Robbe (Backend)"
$d.Content.Find.Execute("Backend", $true, $false, $false, $false, $false, $true, 1, $false, "Robbe (Backend)", 2)

# 3. Backend update body: "Momenteel bezig met" -> "Game logica moet nog verder uitgewerkt worden. " + line break + "Moet nog gebeuren"
$d.Content.Find.Execute("Momenteel bezig met", $true, $false, $false, $false, $false, $true, 1, $false, "Game logica moet nog verder uitgewerkt worden. ^lMoet nog gebeuren", 2)

# 4. "Hardware" heading -> "Herber (Hardware)"
$d.Content.Find.Execute("Hardware", $true, $false, $false, $false, $false, $true, 1, $false, "Herber (Hardware)", 2)

# 5. "Documentatie" heading -> "Jakob (Documentatie)"
$d.Content.Find.Execute("Documentatie", $true, $false, $false, $false, $false, $true, 1, $false, "Jakob (Documentatie)", 2)

# 6. Documentation update body: "Beginnen met Eind documentatie." -> "Beginnen met documentatie voor Milestone 1."
$d.Content.Find.Execute("Beginnen met Eind documentatie.", $true, $false, $false, $false, $false, $true, 1, $false, "Beginnen met documentatie voor Milestone 1.", 2)

Write-Output "done"
